# Updated symbol list on Sat Dec 17 05:20:42 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # The "Price" column stores numeric-looking values as text. Force the
    # cell to stay text (matches the source data's inline-string typing)
    # instead of letting COM auto-coerce a numeric-looking string to a
    # real number.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "231.09"
Set-TextValue "D3"  "22.52"
Set-TextValue "D4"  "5.268"
Set-TextValue "D5"  "0.05548"
Set-TextValue "D6"  "3.379"
Set-TextValue "D8"  "1.058"
Set-TextValue "D9"  "0.7818"
Set-TextValue "D10" "0.1383"
Set-TextValue "D11" "0.07380"
Set-TextValue "D12" "0.03150"
Set-TextValue "D13" "0.02965"
Set-TextValue "D14" "0.09272"
Set-TextValue "D15" "0.001661"
Set-TextValue "D16" "3.266"
Set-TextValue "D17" "0.04773"
Set-TextValue "D18" "0.0005899"
Set-TextValue "D19" "0.006214"
Set-TextValue "D20" "0.005238"
Set-TextValue "D21" "0.001063"
Set-TextValue "D23" "3.914"

# --- Volume(1h) label fix ---
$ws.Range("E27").Value = "26UpBotsUBXT"

Set-TextValue "D40" "0.04006"
Set-TextValue "D41" "0.007139"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# --- Row 42 / 43 swap: CEJI <-> BKEXToken ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003500"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1038"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.009974"
Set-TextValue "D45" "0.00005439"
Set-TextValue "D48" "0.04042"
